$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add two new bullet items ("Sorteren op deadline" and "Filteren op
#    status") to the "Should have" (numId 2) list, right after the
#    "Stijlvolle en responsieve vormgeving (basis CSS/Bootstrap)" item and
#    before the "Could have:" heading.
# ---------------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd("`r") -eq "Stijlvolle en responsieve vormgeving (basis CSS/Bootstrap)") {
        $target = $para
        break
    }
}

$target.Range.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Range.Text = "Sorteren op deadline"

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Filteren op status"

# ---------------------------------------------------------------------------
# 2) In the "Won't have" (numId 4) list, change the highlight color from
#    yellow to red for "Mobiele app", "Meertaligheid" and
#    "Teamtaken (gedeelde takenlijst)". "Meertaligheid" additionally gains a
#    paragraph-mark run-properties block (<w:pPr><w:rPr>) it did not have
#    before. Because the paragraph-mark's own run formatting can't be
#    reached through the regular Font/HighlightColorIndex members, the
#    paragraphs are rewritten wholesale via Range.InsertXML with the exact
#    OOXML that should result, preserving every attribute that the change
#    does not touch.
# ---------------------------------------------------------------------------

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($paragraph, [string]$xml) {
    $paragraph.Range.InsertXML($xml)
}

$mobiele = $null
$meertaligheid = $null
$teamtaken = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.TrimEnd("`r")
    if ($text -eq "Mobiele app") { $mobiele = $para }
    elseif ($text -eq "Meertaligheid") { $meertaligheid = $para }
    elseif ($text -eq "Teamtaken (gedeelde takenlijst)") { $teamtaken = $para }
}

$mobieleXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="112FE84E" w14:textId="77777777" w:rsidR="00EA5236" w:rsidRPr="00F73EDF" w:rsidRDefault="00EA5236" w:rsidP="00EA5236"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:r w:rsidRPr="00F73EDF"><w:rPr><w:b/><w:bCs/><w:highlight w:val="red"/></w:rPr><w:t>Mobiele app</w:t></w:r></w:p>'
Set-ParagraphXml $mobiele $mobieleXml

$meertaligheidXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="75EFBC4B" w14:textId="77777777" w:rsidR="00EA5236" w:rsidRPr="00EA5236" w:rsidRDefault="00EA5236" w:rsidP="00EA5236"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:r w:rsidRPr="00F73EDF"><w:rPr><w:b/><w:bCs/><w:highlight w:val="red"/></w:rPr><w:t>Meertaligheid</w:t></w:r></w:p>'
Set-ParagraphXml $meertaligheid $meertaligheidXml

$teamtakenXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="09177FC9" w14:textId="77777777" w:rsidR="00EA5236" w:rsidRPr="00F73EDF" w:rsidRDefault="00EA5236" w:rsidP="00EA5236"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:r w:rsidRPr="00F73EDF"><w:rPr><w:b/><w:bCs/><w:highlight w:val="red"/></w:rPr><w:t>Teamtaken (gedeelde takenlijst)</w:t></w:r></w:p>'
Set-ParagraphXml $teamtaken $teamtakenXml

Write-Output "Edit complete."
